# "Again Done some modification"
#
# Employees sheet (sheet1 / tab "Employees"):
#   - B6 ("Maddy ")            -> cleared (row 3 now has no middle-name-col value)
#   - B8 ("Piyush  Kumar")     -> "Piyush  Kumar  New"
#   - B9 ("Dharam Kumar")      -> "Dharam Kumar New"
#   - view: top-left visible cell A4 -> A16, selection B9 -> C9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

# Clear the old "Maddy " entry in B6 entirely (cell becomes empty, not just blank text).
$ws.Range("B6").ClearContents()

# Rename entries - set B9 before B8 so new shared strings are appended in the
# same order ("Dharam Kumar New" then "Piyush  Kumar  New").
$ws.Range("B9").Value = "Dharam Kumar New"
$ws.Range("B8").Value = "Piyush  Kumar  New"

# Update the window/sheet view: scroll so row 16 is the top-left visible row,
# and move the selection from B9 to C9.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C9").Select()
